# Add a new "Dozentenumfrage Aufwandsbereitschaft" column (K) to the
# Definition of Ready sheet, mirroring the existing B:J columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Definition of Ready")

# Header for the new column
$ws.Range("K1").Value = "Dozentenumfrage Aufwandsbereitschaft"

# Mark rows 2-6 with an X, matching the existing pattern columns B:J
$ws.Range("K2").Value = "X"
$ws.Range("K3").Value = "X"
$ws.Range("K4").Value = "X"
$ws.Range("K5").Value = "X"
$ws.Range("K6").Value = "X"

# Responsible persons for row 8 (mirrors format of the neighbouring J8 cell)
$ws.Range("K8").Value = "Pütter / Hesse"
$ws.Range("J8").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$excel.CutCopyMode = 0
